$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "Mil-Max sockets" row content (formula in B21, text in D21)
$ws.Range("B21:D21").ClearContents()

# Update the selection to match the final state
$ws.Range("B21:D21").Select()
